$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two forecast entries entirely (not enough history for a
# naive forecast on these rows) - clears C2, E2, C3.
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Remaining values are corrected floating point results from the bugfix in
# the naive component forecaster.
$ws.Range("E3").Value = 6.772115316529947
$ws.Range("C4").Value = -7.266312015249799
$ws.Range("C6").Value = 9.469137444079955
$ws.Range("C7").Value = 3.358206407534969
$ws.Range("E7").Value = 5.745831525574463
$ws.Range("C9").Value = 3.901355411819685
$ws.Range("E9").Value = 3.690459963535031
$ws.Range("C12").Value = 5.246209615995689
$ws.Range("C13").Value = 4.862559663742938
$ws.Range("E13").Value = 4.636196713604379
$ws.Range("C14").Value = 2.76474001115945
$ws.Range("C15").Value = -7.260793671746447
$ws.Range("C16").Value = 4.097586525396246
$ws.Range("C17").Value = 7.824284864703768
$ws.Range("C18").Value = -1.245022353133318
